$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.028.53'
$ws.Range("E2").Value = '  -0.81%  '

$ws.Range("D3").Value = '2.358.62'
$ws.Range("E3").Value = '  -0.33%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.94'
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.19'
$ws.Range("E7").Value = '  +0.88%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  +7.55%  '

$ws.Range("E10").Value = '  -1.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.29'
$ws.Range("E11").Value = '  -0.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.20'
$ws.Range("E12").Value = '  +7.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.26'
$ws.Range("E13").Value = '  +7.01%  '

$ws.Range("E14").Value = '  +1.36%  '

$ws.Range("D15").Value = '2.709.61'
$ws.Range("E15").Value = '  -0.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.54'
$ws.Range("E16").Value = '  -2.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.902'
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").Value = '2.359.13'
$ws.Range("E18").Value = '  -0.36%  '

$ws.Range("D19").Value = '43.943.29'
$ws.Range("E19").Value = '  -1.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.02'
$ws.Range("E20").Value = '  +8.28%  '

$ws.Range("E21").Value = '  -0.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.34'
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '258.78'
$ws.Range("E23").Value = '  +1.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.97'
$ws.Range("E24").Value = '  +22.73%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").Value = '  -4.86%  '

$ws.Range("E27").Value = '  -0.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.79'
$ws.Range("E28").Value = '  +3.74%  '

$ws.Range("E29").Value = '  +1.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.83'
$ws.Range("E30").Value = '  +1.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.43'
$ws.Range("E31").Value = '  +0.89%  '

$ws.Range("E32").Value = '  -1.84%  '

$ws.Range("E33").Value = '  +2.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0759'
$ws.Range("E34").Value = '  +2.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.60'
$ws.Range("E35").Value = '  +7.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.22'
$ws.Range("E36").Value = '  +0.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.76'
$ws.Range("E37").Value = '  -3.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.36'
$ws.Range("E38").Value = '  -1.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.35'
$ws.Range("E39").Value = '  -3.29%  '

$ws.Range("E40").Value = '  +2.55%  '

$ws.Range("E41").Value = '  +13.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.203'
$ws.Range("E42").Value = '  +9.76%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.00'
$ws.Range("E43").Value = '  +1.76%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.89'
$ws.Range("E44").Value = '  -5.48%  '

$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("B46").Value = 'MultiversX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.51'
$ws.Range("E46").Value = '  +13.20%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.75'
$ws.Range("E47").Value = '  +5.84%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.49'
$ws.Range("E48").Value = '  +6.02%  '

$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.14'
$ws.Range("E50").Value = '  +2.38%  '

$ws.Range("E51").Value = '  -0.11%  '
